$d = $word.ActiveDocument

# 1) Insert the VML "J" + underlined "MS" text-box shape as a new leading run
#    in the paragraph that contains "Alex".
$xmlAlexPara = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w10="urn:schemas-microsoft-com:office:word"><w:body><w:p><w:r><w:rPr><w:noProof/><w:lang w:val="en-US" w:eastAsia="zh-TW"/></w:rPr><w:pict><v:shapetype id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="_x0000_s1028" type="#_x0000_t202" style="position:absolute;margin-left:272.4pt;margin-top:13.5pt;width:152.65pt;height:144.55pt;z-index:251660288;mso-width-relative:margin;mso-height-relative:margin" fillcolor="black [3200]" strokecolor="#f2f2f2 [3041]" strokeweight=".25pt"><v:shadow on="t" type="perspective" color="#7f7f7f [1601]" opacity=".5" offset="1pt" offset2="-1pt"/><v:textbox><w:txbxContent><w:p><w:pPr><w:spacing w:before="400" w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Baskerville Old Face" w:hAnsi="Baskerville Old Face"/><w:b/><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr><w:t>J</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Baskerville Old Face" w:hAnsi="Baskerville Old Face"/><w:b/><w:sz w:val="144"/><w:szCs w:val="144"/><w:u w:val="single"/></w:rPr><w:t>MS</w:t></w:r></w:p></w:txbxContent></v:textbox></v:shape></w:pict></w:r><w:r><w:t>Alex</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(2).Range.InsertXML($xmlAlexPara)

# 2) Replace the 3rd of the four blank paragraphs (paragraph 9) with the
#    underlined "MS" (Lucida Fax, 144 half-points) paragraph.
$xmlMsPara = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w10="urn:schemas-microsoft-com:office:word"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Lucida Fax" w:hAnsi="Lucida Fax"/><w:sz w:val="144"/><w:szCs w:val="144"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Fax" w:hAnsi="Lucida Fax"/><w:sz w:val="144"/><w:szCs w:val="144"/><w:u w:val="single"/></w:rPr><w:t>MS</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(9).Range.InsertXML($xmlMsPara)

# 3) Add the (now-referenced) Balloon Text / Balloon Text Char style pair.
$s = $d.Styles.Add("Balloon Text", 1)
$s.ParagraphFormat.SpaceAfter = 0
$s.ParagraphFormat.LineSpacingRule = 0
$s.BaseStyle = "Normal"
$s.Priority = 99
$s.UnhideWhenUsed = $true
$s.Font.Name = "Tahoma"
$s.Font.NameBi = "Tahoma"
$s.Font.Size = 8
$s.Font.SizeBi = 8

$cs = $d.Styles.Add("Balloon Text Char", 2)
$cs.BaseStyle = "DefaultParagraphFont"
$cs.Priority = 99
$cs.Font.Name = "Tahoma"
$cs.Font.NameBi = "Tahoma"
$cs.Font.Size = 8
$cs.Font.SizeBi = 8

$s.LinkStyle = "BalloonTextChar"
$cs.LinkStyle = "BalloonText"

Write-Output "done"
